# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (fund holdings detail) right before the
#    "总计" (Total) summary sheet, cloned from the "2021-Q4" sheet so it
#    inherits the same column layout / styles, then trimmed down to a
#    single data row and repopulated with the 2022-Q1 figures.
# 2. Update the "总计" sheet: add a new first data row for "2022-Q1"
#    (holding count 1, market value 0.01), pushing the existing quarters
#    down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: build the "2022-Q1" worksheet
# ---------------------------------------------------------------------
$source = $wb.Worksheets.Item("2021-Q4")
$total  = $wb.Worksheets.Item("总计")

# Places the clone immediately before "总计". Re-resolve the collection
# afterwards (object refs captured before a sheet insertion can go stale)
# and grab the new sheet by its position, one slot before "总计".
$source.Copy($total, $null)
$totalIndex = $wb.Worksheets.Item("总计").Index
$newSheet = $wb.Worksheets.Item($totalIndex - 1)

# The source sheet had 4 data rows (rows 2-5); 2022-Q1 only needs 1, so
# drop the extra three and let the remaining row get overwritten below.
$newSheet.Rows.Item(3).Delete()
$newSheet.Rows.Item(3).Delete()
$newSheet.Rows.Item(3).Delete()

# Force text storage for the numeric-looking identifiers/figures (matches
# the source data, which stores these columns as text) while leaving the
# cell format otherwise untouched.
$newSheet.Range("B2").Value = "'003981"
$newSheet.Range("C2").Value = "中银证券瑞益灵活配置混合C"
$newSheet.Range("D2").Value = "'0.21"
$newSheet.Range("E2").Value = "'89.21"
$newSheet.Range("F2").Value = "'4.39"
$newSheet.Range("G2").Value = "'0.0092"
$newSheet.Range("H2").Value2 = 4

$newSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------
# Step 2: update the "总计" worksheet with the new 2022-Q1 row
# ---------------------------------------------------------------------
# Re-resolve "总计" again: the rename in step 1 is one more structural
# change that can stale out previously captured references.
$total = $wb.Worksheets.Item("总计")

# Shift the existing B:D data down by one row (from the bottom up so we
# never clobber a row before it has been copied).
for ($r = 7; $r -ge 3; $r--) {
    $srcRow = $r - 1
    $total.Range("B$r").Value2 = $total.Range("B$srcRow").Value2
    $total.Range("C$r").Value2 = $total.Range("C$srcRow").Value2
    $total.Range("D$r").Value2 = $total.Range("D$srcRow").Value2
}

# Carry the A-column (index 0..5) style down to the newly-used row 7
# before refreshing all the index values.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122, $null, $false, $false)

# New first data row: 2022-Q1, 1 holding, 0.01 billion market value.
$total.Range("B2").Value2 = "2022-Q1"
$total.Range("C2").Value2 = 1
$total.Range("D2").Value2 = 0.01

# Refresh the sequential row index column (A) for every data row.
for ($r = 2; $r -le 7; $r++) {
    $total.Range("A$r").Value2 = $r - 2
}
